$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.203.65"
$ws.Range("E2").Value = "  +0.60%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.803.82"
$ws.Range("E3").Value = "  -0.05%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.45"
$ws.Range("E5").Value = "  +0.73%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.48"
$ws.Range("E6").Value = "  -1.08%  "

# Row 7
$ws.Range("E7").Value = "  -0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  -0.41%  "

# Row 9
$ws.Range("E9").Value = "  -0.96%  "

# Row 10
$ws.Range("E10").Value = "  +0.63%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.48"
$ws.Range("E11").Value = "  +2.77%  "

# Row 12
$ws.Range("E12").Value = "  -1.39%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "35.93"
$ws.Range("E13").Value = "  -0.53%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.440.20"
$ws.Range("E14").Value = "  -0.06%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.859.71"
$ws.Range("E15").Value = "  +2.20%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "68.186.66"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.47"
$ws.Range("E17").Value = "  -0.79%  "

# Row 18
$ws.Range("E18").Value = "  +2.23%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.11"
$ws.Range("E19").Value = "  -0.09%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "462.81"
$ws.Range("E20").Value = "  +0.27%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.74"
$ws.Range("E21").Value = "  -1.96%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.702"
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000150"
$ws.Range("E23").Value = "  -3.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.09"
$ws.Range("E24").Value = "  -0.66%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.08"
$ws.Range("E25").Value = "  -0.15%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.12"
$ws.Range("E26").Value = "  +0.45%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.03"
$ws.Range("E27").Value = "  +0.16%  "

# Row 28
$ws.Range("E28").Value = "  -0.14%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.952.14"
$ws.Range("E29").Value = "  -0.05%  "

# Row 30
$ws.Range("E30").Value = "  -4.96%  "

# Row 31
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.24"
$ws.Range("E31").Value = "  -0.87%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.36"
$ws.Range("E32").Value = "  +1.32%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.45"
$ws.Range("E33").Value = "  -1.10%  "

# Row 34
$ws.Range("E34").Value = "  +0.02%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.06"
$ws.Range("E35").Value = "  -0.38%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.100"
$ws.Range("E36").Value = "  -0.08%  "

# Row 37
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.33"
$ws.Range("E37").Value = "  -3.01%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.139"
$ws.Range("E38").Value = "  +0.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.83"
$ws.Range("E39").Value = "  +0.64%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.990"
$ws.Range("E40").Value = "  -0.43%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  -0.06%  "

# Row 43
$ws.Range("E43").Value = "  +0.64%  "

# Row 44
$ws.Range("E44").Value = "  -1.42%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.13"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "151.61"
$ws.Range("E46").Value = "  +1.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.39"

# Row 48
$ws.Range("E48").Value = "  +2.32%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "397.02"
$ws.Range("E49").Value = "  -0.17%  "

# Row 50
$ws.Range("E50").Value = "  +4.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.67"
$ws.Range("E51").Value = "  +0.48%  "
